# Update odds values in row 3 and row 4 of Sheet1 to reflect the latest
# FlashScore odds snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 updates ---
$ws.Range("K3").Value  = 2.75
$ws.Range("M3").Value  = 1.02
$ws.Range("N3").Value  = 19
$ws.Range("O3").Value  = 1.13
$ws.Range("P3").Value  = 6
$ws.Range("Q3").Value  = 1.48
$ws.Range("R3").Value  = 2.6
$ws.Range("U3").Value  = 1.8
$ws.Range("V3").Value  = 1.91
$ws.Range("W3").Value  = 9.5
$ws.Range("X3").Value  = 7.5
$ws.Range("Y3").Value  = 9
$ws.Range("Z3").Value  = 9
$ws.Range("AC3").Value = 19
$ws.Range("AE3").Value = 19
$ws.Range("AG3").Value = 26
$ws.Range("AH3").Value = 51
$ws.Range("AM3").Value = 201
$ws.Range("AN3").Value = 3.5
$ws.Range("AU3").Value = 8.5
$ws.Range("AY3").Value = 34
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 201

# --- Row 4 updates ---
$ws.Range("G4").Value  = 1.8
$ws.Range("H4").Value  = 3.3
$ws.Range("J4").Value  = 2.5
$ws.Range("O4").Value  = 1.53
$ws.Range("P4").Value  = 2.5
$ws.Range("Q4").Value  = 2.6
$ws.Range("R4").Value  = 1.48
$ws.Range("AC4").Value = 6.5
$ws.Range("AN4").Value = 3.5
